# feat: add 2022-Q1 data
#
# - Insert a new worksheet "2022-Q1" positioned right before the "总计"
#   (Total) sheet, with the same column layout as the other quarterly
#   sheets (基金代码/基金名称/基金规模/股票总仓位/仓位占比/持有市值(亿元)/仓位排名),
#   populated with the Q1-2022 fund-holding rows.
# - Insert a new top row into the "总计" sheet summarizing 2022-Q1
#   (6 holdings, 0.75 亿元), pushing the existing quarters down by one row.

$wb = $excel.ActiveWorkbook

$totalSheetRef = $wb.Worksheets.Item("总计")
$q4Sheet = $wb.Worksheets.Item("2021-Q4")

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q1" sheet right before "总计".
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($totalSheetRef)
$newSheet.Name = "2022-Q1"

# Re-fetch the "总计" sheet by name: adding a sheet can shift/rebind
# earlier object references, so grab a fresh handle to be safe.
$totalSheet = $wb.Worksheets.Item("总计")

# Carry over the header-row (B1:H1) and index-column (A2:A7) cell
# formatting from the "2021-Q4" sheet, which already uses the shared
# bold/bordered "s=2" style used throughout this workbook.
$q4Sheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)   # xlPasteFormats
$q4Sheet.Range("A2:A5").Copy()
$newSheet.Range("A2:A7").PasteSpecial(-4122)   # xlPasteFormats

function Set-TextCell($sheet, $addr, $text) {
    # Force the value to be written as a genuine text string (matching
    # the workbook's t="inlineStr" convention) instead of letting COM
    # auto-coerce numeric-looking strings ("003721", "6.31", ...) into
    # numbers. Writing a text formula and then collapsing it down to a
    # static value keeps the cell's existing number format/style intact.
    $cell = $sheet.Range($addr)
    $escaped = $text -replace '"', '""'
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)   # xlPasteValues
}

# Headers
Set-TextCell $newSheet "B1" "基金代码"
Set-TextCell $newSheet "C1" "基金名称"
Set-TextCell $newSheet "D1" "基金规模"
Set-TextCell $newSheet "E1" "股票总仓位"
Set-TextCell $newSheet "F1" "仓位占比"
Set-TextCell $newSheet "G1" "持有市值(亿元)"
Set-TextCell $newSheet "H1" "仓位排名"

$q1_2022 = @(
    @("003721", "易方达标普信息科技指数（QDII-LOF）美元", "6.31", "93.58", "3.21", "0.2026", 4),
    @("161128", "易方达标普信息科技指数（QDII-LOF）人民币", "6.31", "93.58", "3.21", "0.2026", 4),
    @("000043", "嘉实美国成长股票(QDII) -人民币", "14.64", "94.24", "1.15", "0.1684", 10),
    @("000044", "嘉实美国成长股票(QDII) - 美元现汇", "14.64", "94.24", "1.15", "0.1684", 10),
    @("519981", "长信美国标准普尔100等权重指数增强(QDII)", "0.47", "84.16", "0.87", "0.0041", 7),
    @("011706", "长信美国标准普尔100等权重指数增强(QDII) - 美元", "0.47", "84.16", "0.87", "0.0041", 7)
)

$r = 2
foreach ($row in $q1_2022) {
    $newSheet.Range("A$r").Value2 = $r - 2
    Set-TextCell $newSheet "B$r" $row[0]
    Set-TextCell $newSheet "C$r" $row[1]
    Set-TextCell $newSheet "D$r" $row[2]
    Set-TextCell $newSheet "E$r" $row[3]
    Set-TextCell $newSheet "F$r" $row[4]
    Set-TextCell $newSheet "G$r" $row[5]
    $newSheet.Range("H$r").Value2 = $row[6]
    $r++
}

# ---------------------------------------------------------------------
# 2. Insert the 2022-Q1 summary row at the top of the "总计" sheet,
#    shifting the existing quarters down by one row.
# ---------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()

# Inserting a row clones the row-above's formatting onto the new cells
# (bordered "s=3" style leaking onto B2:D2). Reset B2:D2 to the plain,
# unstyled look the other data rows use, and copy the bold/bordered
# index-column style from A3 (the row that held this formatting before
# the insert) onto the new A2.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)   # xlPasteFormats
$totalSheet.Range("B2:D2").ClearFormats()

$totalSheet.Range("A2").Value2 = 0
Set-TextCell $totalSheet "B2" "2022-Q1"
$totalSheet.Range("C2").Value2 = 6
$totalSheet.Range("D2").Value2 = 0.75

# Renumber the index column (A) sequentially for the rows that shifted
# down (they keep their original 0-based value otherwise).
$totalSheet.Range("A3").Value2 = 1
$totalSheet.Range("A4").Value2 = 2
$totalSheet.Range("A5").Value2 = 3
$totalSheet.Range("A6").Value2 = 4

Write-Host "2022-Q1 sheet added and 总计 updated"
